# Fix problem with email specifications in cat
# - Correct the "skeleton" target paths on the "cat" sheet so the target
#   directories actually match the copied source sub-paths (fonts, js/ckeditor,
#   focussearch) and the email target now points at "/email" instead of "/".
# - Restore "-Instructions" as the active/selected tab (it had drifted to "cat").

$wb = $excel.ActiveWorkbook

$catSheet = $wb.Worksheets.Item("cat")

$catSheet.Range("C10").Value = "/web/<progDir>/fonts"
$catSheet.Range("C11").Value = "/web/<progDir>/js/ckeditor"
$catSheet.Range("C12").Value = "/web/<progDir>/focussearch"
$catSheet.Range("C13").Value = "/email"

# Update the remembered selection on the "cat" sheet (this also momentarily
# activates it), then activate "-Instructions" so it becomes the active tab
# again, leaving "cat" unselected.
[void]$catSheet.Range("A10").Select()

$instructionsSheet = $wb.Worksheets.Item("-Instructions")
[void]$instructionsSheet.Activate()
[void]$instructionsSheet.Range("C17").Select()
